$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 25 de Junio de 2020 a las 00:04'

# Row 4
$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 2457574
$ws.Cells.Item(4, 3).Value = 33406
$ws.Cells.Item(4, 4).Value = 1033275
$ws.Cells.Item(4, 5).Value = 1300113
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 713
$ws.Cells.Item(4, 8).Value = 124186

# Row 5
$ws.Cells.Item(5, 1).Value = 'Brasil'
$ws.Cells.Item(5, 2).Value = 1188631
$ws.Cells.Item(5, 3).Value = 37152
$ws.Cells.Item(5, 4).Value = 649908
$ws.Cells.Item(5, 5).Value = 484893
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 1059
$ws.Cells.Item(5, 8).Value = 53830

# Row 14
$ws.Cells.Item(14, 1).Value = 'Alemania'
$ws.Cells.Item(14, 2).Value = 193217
$ws.Cells.Item(14, 3).Value = 439
$ws.Cells.Item(14, 4).Value = 176300
$ws.Cells.Item(14, 5).Value = 7914
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 17
$ws.Cells.Item(14, 8).Value = 9003

# Row 21
$ws.Cells.Item(21, 1).Value = 'Sudafrica'
$ws.Cells.Item(21, 2).Value = 111796
$ws.Cells.Item(21, 3).Value = 5688
$ws.Cells.Item(21, 4).Value = 56874
$ws.Cells.Item(21, 5).Value = 52717
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 103
$ws.Cells.Item(21, 8).Value = 2205

# Row 22
$ws.Cells.Item(22, 1).Value = 'Canada'
$ws.Cells.Item(22, 2).Value = 102228
$ws.Cells.Item(22, 3).Value = 265
$ws.Cells.Item(22, 4).Value = 64976
$ws.Cells.Item(22, 5).Value = 28769
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 29
$ws.Cells.Item(22, 8).Value = 8483

# Row 50
$ws.Cells.Item(50, 1).Value = 'Barein'
$ws.Cells.Item(50, 2).Value = 23570
$ws.Cells.Item(50, 3).Value = 508
$ws.Cells.Item(50, 4).Value = 17977
$ws.Cells.Item(50, 5).Value = 5524
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 2
$ws.Cells.Item(50, 8).Value = 69

# Row 51
$ws.Cells.Item(51, 1).Value = 'Israel'
$ws.Cells.Item(51, 2).Value = 22044
$ws.Cells.Item(51, 3).Value = 532
$ws.Cells.Item(51, 4).Value = 15940
$ws.Cells.Item(51, 5).Value = 5796
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 308

# Row 64
$ws.Cells.Item(64, 1).Value = 'Camerun'
$ws.Cells.Item(64, 2).Value = 12592
$ws.Cells.Item(64, 3).Value = 322
$ws.Cells.Item(64, 4).Value = 10100
$ws.Cells.Item(64, 5).Value = 2179
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 313

# Row 65
$ws.Cells.Item(65, 1).Value = 'Corea del Sur'
$ws.Cells.Item(65, 2).Value = 12535
$ws.Cells.Item(65, 3).Value = 51
$ws.Cells.Item(65, 4).Value = 10930
$ws.Cells.Item(65, 5).Value = 1324
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 281

# Row 73
$ws.Cells.Item(73, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(73, 2).Value = 8164
$ws.Cells.Item(73, 3).Value = 260
$ws.Cells.Item(73, 4).Value = 3419
$ws.Cells.Item(73, 5).Value = 4687
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 58

# Row 88
$ws.Cells.Item(88, 1).Value = 'Bulgaria'
$ws.Cells.Item(88, 2).Value = 4242
$ws.Cells.Item(88, 3).Value = 128
$ws.Cells.Item(88, 4).Value = 2263
$ws.Cells.Item(88, 5).Value = 1770
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(88, 8).Value = 209

# Row 89
$ws.Cells.Item(89, 1).Value = 'Venezuela'
$ws.Cells.Item(89, 2).Value = 4187
$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(89, 4).Value = 1327
$ws.Cells.Item(89, 5).Value = 2825
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 35

# Row 90
$ws.Cells.Item(90, 1).Value = 'Luxemburgo'
$ws.Cells.Item(90, 2).Value = 4140
$ws.Cells.Item(90, 3).Value = 7
$ws.Cells.Item(90, 4).Value = 3965
$ws.Cells.Item(90, 5).Value = 65
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 110

# Row 110
$ws.Cells.Item(110, 1).Value = 'Sudan del Sur'
$ws.Cells.Item(110, 2).Value = 1942
$ws.Cells.Item(110, 3).Value = 12
$ws.Cells.Item(110, 4).Value = 224
$ws.Cells.Item(110, 5).Value = 1682
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 36

# Row 144
$ws.Cells.Item(144, 1).Value = 'Santo Tome y Principe'
$ws.Cells.Item(144, 2).Value = 710
$ws.Cells.Item(144, 3).Value = 3
$ws.Cells.Item(144, 4).Value = 211
$ws.Cells.Item(144, 5).Value = 486
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 1
$ws.Cells.Item(144, 8).Value = 13

# Row 152
$ws.Cells.Item(152, 1).Value = 'Zimbabue'
$ws.Cells.Item(152, 2).Value = 530
$ws.Cells.Item(152, 3).Value = 5
$ws.Cells.Item(152, 4).Value = 123
$ws.Cells.Item(152, 5).Value = 401
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 6

# Row 168
$ws.Cells.Item(168, 1).Value = 'Islas Caimanes'
$ws.Cells.Item(168, 2).Value = 196
$ws.Cells.Item(168, 3).Value = 1
$ws.Cells.Item(168, 4).Value = 169
$ws.Cells.Item(168, 5).Value = 26
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 1

# Row 182
$ws.Cells.Item(182, 1).Value = 'Botsuana'
$ws.Cells.Item(182, 2).Value = 92
$ws.Cells.Item(182, 3).Value = 3
$ws.Cells.Item(182, 4).Value = 25
$ws.Cells.Item(182, 5).Value = 66
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 1

# Row 202
$ws.Cells.Item(202, 1).Value = 'Dominica'
$ws.Cells.Item(202, 2).Value = 18
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 18
$ws.Cells.Item(202, 5).Value = 0
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

# Row 203
$ws.Cells.Item(203, 1).Value = 'Fiyi'
$ws.Cells.Item(203, 2).Value = 18
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(203, 4).Value = 18
$ws.Cells.Item(203, 5).Value = 0
$ws.Cells.Item(203, 6).Value = 0
$ws.Cells.Item(203, 7).Value = 0
$ws.Cells.Item(203, 8).Value = 0

# Row 208
$ws.Cells.Item(208, 1).Value = 'Groenlandia'
$ws.Cells.Item(208, 2).Value = 13
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 13
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

# Row 209
$ws.Cells.Item(209, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(209, 2).Value = 13
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 13
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0
